# Applies the "6.0.0" release update to the StructureDefinition workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" -------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicate "Contact" / "No display for ContactDetail" row;
# it becomes the new "Jurisdiction" / "United States of America" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" row and is removed entirely,
# shifting everything below it up by one row.
$ws.Rows.Item(11).Delete()

# --- Sheet "Elements" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# The root Extension element's Short/Definition now mirror the new
# Title / Description of the structure definition.
$ws2.Range("K2").Value = "Episode Admit Count"
$ws2.Range("L2").Value = "Number of admissions related to the episode of care"
